$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 102: add E102 = "No" (new cell added between existing D102 and F102) ---
$ws.Cells.Item(102, 5).Value = "No"

# --- Row 109: new row ---
# A109 is a literal text string (trailing backtick makes it not a real date), so it must
# stay text rather than being interpreted/parsed as a date value.
$ws.Cells.Item(109, 1).Value = "3/5/2025``"
$ws.Cells.Item(109, 2).Value = "Amazon"
$ws.Cells.Item(109, 3).Value = "sr. data scientist"
$ws.Cells.Item(109, 4).Value = "Interview w/ management"
$ws.Cells.Item(109, 5).Value = "No"
$ws.Cells.Item(109, 6).Value = "the energy job, same url as before"

# --- Row 110: new row ---
# A110 is a real date (3/7/2025 => serial 45723). Assign the numeric serial directly so the
# engine doesn't auto-generate a brand-new custom number format, then apply the same
# "m/d/yy" format used elsewhere in column A, which maps back onto the workbook's existing
# built-in date style instead of creating a new one.
$ws.Cells.Item(110, 1).Value = 45723
$ws.Cells.Item(110, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(110, 2).Value = "DTN"
$ws.Cells.Item(110, 3).Value = "Sr Data Scientist"
$ws.Cells.Item(110, 4).Value = "First interview"
$ws.Cells.Item(110, 6).Value = "same url as before"

# --- Row 111: new row ---
$ws.Cells.Item(111, 1).Value = 45724
$ws.Cells.Item(111, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(111, 2).Value = "MISO"
$ws.Cells.Item(111, 3).Value = "Senior Engineer Market Evaluation"
$ws.Cells.Item(111, 4).Value = "application, MISO market efficiency evaluator"
$ws.Cells.Item(111, 6).Value = "https://recruiting.ultipro.com/MID1029MISO/JobBoard/362b6b1d-f1c3-46f5-9554-4aa90e2bda64/OpportunityDetail?opportunityId=ae933cb8-99a0-42af-ae01-b41431abde3a"

Write-Host "edit complete"
